# Vuvlo_calc.xlsx update:
#  - finalize CAM
#  - add classic drawing for TOP
#  - recalculate VUVLO for new 1.2 schematics

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("PowerConditioner")
$ws2 = $wb.Worksheets.Item("9V")

# ---------------------------------------------------------------------------
# Sheet "PowerConditioner" (sheet1)
# ---------------------------------------------------------------------------

# --- Fully clear (content + formatting) the whole block that is being
# restructured, so stale formulas/values/styles from the old 3s/4s/5s table
# don't leave any residue, then rebuild it from scratch.
$ws1.Range("A9:K19").Clear()
$ws1.Range("K1:R3").Clear()

# --- Row 5/6: new divider values for the 1.2V reference design
$ws1.Range("B5").Value = 400
$ws1.Range("B6").Value = 1.2

# --- Row 9 header row (Vuvlo table)
$ws1.Range("A9").Value = "Vbg = 1,21"
$ws1.Range("C9").Value = "Target stop V"
$ws1.Range("D9").Value = "R2 [kOhm]"
$ws1.Range("E9").Value = "R2 chosen"
$ws1.Range("F9").Value = "Stop voltage"
$ws1.Range("G9").Value = "stop voltage/cell"

# --- Row 10: single remaining "5s" Vuvlo calculation
$ws1.Range("A10").Value = "Vuvlo5s = 15   "
$ws1.Range("B10").Value = 5
$ws1.Range("C10").Value = 15
$ws1.Range("D10").NumberFormat = "0.0"
$ws1.Range("D10").Formula = '=$B$6*(1 + $B$5/C10)'
$ws1.Range("E10").Value = 33
$ws1.Range("F10").NumberFormat = "0.0"
$ws1.Range("F10").Formula = '=$B$6*(1 + $B$5/E10)'
$ws1.Range("G10").Formula = '=F10/B10'

# --- Row 11: leftover styled (blank) cells
$ws1.Range("D11").NumberFormat = "0.0"
$ws1.Range("G11").NumberFormat = "0.0"
$ws1.Range("K11").NumberFormat = "0.0"

# --- Row 12: Cin block + leftover styled cell
$ws1.Range("A12").Value = "Cin"
$ws1.Range("B12").Value = "20u"
$ws1.Range("C12").Value = "35v"
$ws1.Range("K12").NumberFormat = "0.0"

# --- Row 13: Cout block
$ws1.Range("A13").Value = "Cout"
$ws1.Range("B13").Value = "100n"
$ws1.Range("C13").Value = "35v"

# --- Row 15/16: SETI / R2 chosen table
$ws1.Range("A15").Value = "SETI"
$ws1.Range("B15").Value = "4.2A (max)"
$ws1.Range("C15").Value = "3A"
$ws1.Range("D15").Value = "3.3A"
$ws1.Range("E15").Value = "2.3A"

$ws1.Range("B16").Value = "2.4kOhm"
$ws1.Range("C16").Value = "3.6kOhm"
$ws1.Range("D16").Value = "3.3kOhm"
$ws1.Range("E16").Value = "4.7kOhm"

# --- New OVLO block (K1:O3, R3)
$ws1.Range("K1").Value = "OVLO"

$ws1.Range("L2").Value = "Stop V"
$ws1.Range("M2").Value = "R2 [kOhm]"
$ws1.Range("N2").Value = "R2 chosen"
$ws1.Range("O2").Value = "Stop voltage"

$ws1.Range("K3").Value = "Vovlo = 24   "
$ws1.Range("L3").Value = 24
$ws1.Range("M3").NumberFormat = "0.0"
$ws1.Range("M3").Formula = '=$B$6*(1 + $B$5/L3)'
$ws1.Range("N3").NumberFormat = "0.0"
$ws1.Range("N3").Value = 20
$ws1.Range("O3").NumberFormat = "0.0"
$ws1.Range("O3").Formula = '=$B$6*(1 + $B$5/N3)'
$ws1.Range("R3").Formula = '=N3/#REF!'

# --- Selection marker matches the saved state of the authored workbook
$ws1.Range("E20").Select()

# ---------------------------------------------------------------------------
# Sheet "9V" (sheet2) - shared-string table renumbered only, no content change
# ---------------------------------------------------------------------------
# (left untouched - text values for this sheet are unchanged by the edit)
